$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new rows into the "Hab populations" table (sorted by
#    population, descending). The table is sorted, so the two new entries
#    ("Lockyer" pop 70000 and "Torre Verde" pop 20000) need to land in the
#    correct spots:
#      - Lockyer (70000) goes between Argoed (80000, row 20) and
#        Fields of Steel (60000, old row 21) -> new row 21
#      - Torre Verde (20000) goes between Ashoka (25000, old row 25) and
#        Korolev (10000, old row 26) -> new row 27 (after the first insert
#        shifted everything below row 21 down by one)
# ---------------------------------------------------------------------------
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(27).Insert()

# New row 21: Lockyer
$ws.Range("A21").Value() = "Lockyer"
$ws.Range("B21").Value() = "Hellas Plaintia"
$ws.Range("D21").Value() = "EP Meltwater"
$ws.Range("E21").Value() = "HEL 7"
$ws.Range("F21").Value() = "Crater"
$ws.Range("G21").Value() = 70000
$ws.Range("H21").Value() = 0.3
$ws.Range("I21").Formula = "=G21*H21"

# New row 27: Torre Verde
$ws.Range("A27").Value() = "Torre Verde"
$ws.Range("B27").Value() = "Orbit"
$ws.Range("C27").Value() = "PIA"
$ws.Range("D27").Value() = "Mars.doc"
$ws.Range("E27").Value() = "ORB 2"
$ws.Range("F27").Value() = "Torus"
$ws.Range("G27").Value() = 20000
$ws.Range("H27").Value() = 0.1
$ws.Range("I27").Formula = "=G27*H27"
$ws.Range("J27").Value() = "Party Zone for orbital workers"

# ---------------------------------------------------------------------------
# 2. Update a handful of existing entries (now shifted down by two rows
#    because of the inserts above).
# ---------------------------------------------------------------------------

# "Lu Xing" (now row 33): orbital torus re-purposed into a cluster
$ws.Range("F33").Value() = "Cluster"
$ws.Range("H33").Value() = 0.1
$ws.Range("J33").Value() = "Biological production cluster"

# "Durango" (now row 40): comment corrected ("Faa Jing" -> "Fa Jing")
$ws.Range("J40").Value() = "Taken over by Fa Jing"

# "Ptah" (now row 42): citizen percentage revised down
$ws.Range("H42").Value() = 0.2

# "Piros Lyuk" (now row 45): location name completed
$ws.Range("B45").Value() = "S/Mare Australe"

# ---------------------------------------------------------------------------
# 3. The row inserts above auto-adjusted most formulas, but a couple of
#    "overflow" helper ranges are intentionally left untouched by the
#    author, and the senate COUNTIF helper formulas need their starting
#    bound corrected manually. Force the exact formula text.
# ---------------------------------------------------------------------------
$ws.Range("L5").Formula = "=SUM(G1:G124)"
$ws.Range("L9").Formula = "=SUM(I3:I125)"

$ws.Range("N32").Formula = "=COUNTIF(D15:D1014,""*D*"")+COUNTIF(D15:D1014,""DD*"")+COUNTIF(D15:D1014,""DDD"")"
$ws.Range("N33").Formula = "=COUNTIF(D15:D114,""*X*"")+COUNTIF(D15:D114,""XX*"")+COUNTIF(D15:D1014,""XXX"")"
$ws.Range("N35").Formula = "=COUNTIF(D15:D114,""*F*"")+COUNTIF(D15:D114,""FF*"")+COUNTIF(D15:D1014,""FFF"")"
$ws.Range("N36").Formula = "=COUNTIF(D15:D114,""*R*"")+COUNTIF(D15:D114,""RR*"")+COUNTIF(D15:D1014,""RRR"")"
$ws.Range("N41").Formula = "=SUM(N34:N40)"

$wb.Save()
